# Edits for Jogos_da_Semana_FlashScore_2024-11-06.xlsx
# - Remove the two ENGLAND - CHAMPIONSHIP fixtures (previously rows 8 & 9)
# - Re-order the "Odd_CS_*" correct-score columns AG:AM (left-rotate by one,
#   Odd_CS_4-4 moves from AG to AM) on the header and on every remaining data row
# - Refresh match data / odds for the remaining rows (new fixtures slot into
#   rows 2 & 3, rows 4-7 keep their fixtures but get updated odds)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove rows 8 and 9 (Luton-Cardiff, Preston-Sunderland) ---
$ws.Range("A8:A9").EntireRow.Delete()

# --- Row 1 ---
$ws.Cells.Item(1, 33).Value = 'Odd_CS_0-1'  # AG1
$ws.Cells.Item(1, 34).Value = 'Odd_CS_0-2'  # AH1
$ws.Cells.Item(1, 35).Value = 'Odd_CS_1-2'  # AI1
$ws.Cells.Item(1, 36).Value = 'Odd_CS_0-3'  # AJ1
$ws.Cells.Item(1, 37).Value = 'Odd_CS_1-3'  # AK1
$ws.Cells.Item(1, 38).Value = 'Odd_CS_2-3'  # AL1
$ws.Cells.Item(1, 39).Value = 'Odd_CS_4-4'  # AM1

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = 'E9gEzDwN'  # A2
$ws.Cells.Item(2, 3).Value = '19:30'  # C2
$ws.Cells.Item(2, 5).Value = 'Boca Juniors'  # E2
$ws.Cells.Item(2, 6).Value = 'Godoy Cruz'  # F2
$ws.Cells.Item(2, 7).Value = 1.73  # G2
$ws.Cells.Item(2, 8).Value = 3.4  # H2
$ws.Cells.Item(2, 9).Value = 5.25  # I2
$ws.Cells.Item(2, 10).Value = 2.4  # J2
$ws.Cells.Item(2, 11).Value = 2  # K2
$ws.Cells.Item(2, 12).Value = 6  # L2
$ws.Cells.Item(2, 13).Value = 1.11  # M2
$ws.Cells.Item(2, 14).Value = 6.5  # N2
$ws.Cells.Item(2, 15).Value = 1.5  # O2
$ws.Cells.Item(2, 16).Value = 2.5  # P2
$ws.Cells.Item(2, 17).Value = 2.5  # Q2
$ws.Cells.Item(2, 18).Value = 1.5  # R2
$ws.Cells.Item(2, 19).Value = 1.53  # S2
$ws.Cells.Item(2, 20).Value = 2.38  # T2
$ws.Cells.Item(2, 21).Value = 2.25  # U2
$ws.Cells.Item(2, 22).Value = 1.57  # V2
$ws.Cells.Item(2, 23).Value = 5  # W2
$ws.Cells.Item(2, 24).Value = 6.5  # X2
$ws.Cells.Item(2, 25).Value = 9  # Y2
$ws.Cells.Item(2, 26).Value = 13  # Z2
$ws.Cells.Item(2, 27).Value = 17  # AA2
$ws.Cells.Item(2, 29).Value = 6.5  # AC2
$ws.Cells.Item(2, 30).Value = 7  # AD2
$ws.Cells.Item(2, 31).Value = 23  # AE2
$ws.Cells.Item(2, 33).Value = 11  # AG2
$ws.Cells.Item(2, 34).Value = 26  # AH2
$ws.Cells.Item(2, 35).Value = 19  # AI2
$ws.Cells.Item(2, 36).Value = 67  # AJ2
$ws.Cells.Item(2, 37).Value = 51  # AK2
$ws.Cells.Item(2, 38).Value = 51  # AL2
$ws.Cells.Item(2, 39).Value = 201  # AM2
$ws.Cells.Item(2, 40).Value = 3.4  # AN2
$ws.Cells.Item(2, 41).Value = 9.5  # AO2
$ws.Cells.Item(2, 42).Value = 26  # AP2
$ws.Cells.Item(2, 43).Value = 34  # AQ2
$ws.Cells.Item(2, 44).Value = 67  # AR2
$ws.Cells.Item(2, 46).Value = 2.38  # AT2
$ws.Cells.Item(2, 47).Value = 10  # AU2
$ws.Cells.Item(2, 48).Value = 81  # AV2
$ws.Cells.Item(2, 49).Value = 7  # AW2
$ws.Cells.Item(2, 50).Value = 34  # AX2
$ws.Cells.Item(2, 51).Value = 41  # AY2
$ws.Cells.Item(2, 52).Value = 126  # AZ2
$ws.Cells.Item(2, 53).Value = 201  # BA2
$ws.Cells.Item(2, 54).Value = 501  # BB2
$ws.Cells.Item(2, 56).Value = 151  # BD2

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value = 'bFzviYN3'  # A3
$ws.Cells.Item(3, 3).Value = '21:30'  # C3
$ws.Cells.Item(3, 5).Value = 'Instituto'  # E3
$ws.Cells.Item(3, 6).Value = 'River Plate'  # F3
$ws.Cells.Item(3, 7).Value = 3.6  # G3
$ws.Cells.Item(3, 8).Value = 3.1  # H3
$ws.Cells.Item(3, 9).Value = 2.2  # I3
$ws.Cells.Item(3, 10).Value = 4.33  # J3
$ws.Cells.Item(3, 11).Value = 1.95  # K3
$ws.Cells.Item(3, 12).Value = 3  # L3
$ws.Cells.Item(3, 13).Value = 1.1  # M3
$ws.Cells.Item(3, 14).Value = 7  # N3
$ws.Cells.Item(3, 21).Value = 2.05  # U3
$ws.Cells.Item(3, 22).Value = 1.7  # V3
$ws.Cells.Item(3, 23).Value = 8.5  # W3
$ws.Cells.Item(3, 24).Value = 17  # X3
$ws.Cells.Item(3, 25).Value = 13  # Y3
$ws.Cells.Item(3, 26).Value = 41  # Z3
$ws.Cells.Item(3, 27).Value = 34  # AA3
$ws.Cells.Item(3, 30).Value = 6  # AD3
$ws.Cells.Item(3, 31).Value = 19  # AE3
$ws.Cells.Item(3, 32).Value = 67  # AF3
$ws.Cells.Item(3, 33).Value = 6  # AG3
$ws.Cells.Item(3, 34).Value = 9  # AH3
$ws.Cells.Item(3, 35).Value = 9.5  # AI3
$ws.Cells.Item(3, 36).Value = 21  # AJ3
$ws.Cells.Item(3, 37).Value = 21  # AK3
$ws.Cells.Item(3, 38).Value = 41  # AL3
$ws.Cells.Item(3, 39).Value = 501  # AM3
$ws.Cells.Item(3, 40).Value = 5.5  # AN3
$ws.Cells.Item(3, 41).Value = 21  # AO3
$ws.Cells.Item(3, 42).Value = 34  # AP3
$ws.Cells.Item(3, 43).Value = 81  # AQ3
$ws.Cells.Item(3, 44).Value = 126  # AR3
$ws.Cells.Item(3, 45).Value = 301  # AS3
$ws.Cells.Item(3, 47).Value = 9  # AU3
$ws.Cells.Item(3, 48).Value = 67  # AV3
$ws.Cells.Item(3, 49).Value = 4  # AW3
$ws.Cells.Item(3, 50).Value = 13  # AX3
$ws.Cells.Item(3, 51).Value = 29  # AY3
$ws.Cells.Item(3, 52).Value = 41  # AZ3
$ws.Cells.Item(3, 53).Value = 81  # BA3
$ws.Cells.Item(3, 54).Value = 251  # BB3

# --- Row 4 ---
$ws.Cells.Item(4, 15).Value = 1.5  # O4
$ws.Cells.Item(4, 16).Value = 2.63  # P4
$ws.Cells.Item(4, 33).Value = 8  # AG4
$ws.Cells.Item(4, 34).Value = 15  # AH4
$ws.Cells.Item(4, 35).Value = 12  # AI4
$ws.Cells.Item(4, 36).Value = 34  # AJ4
$ws.Cells.Item(4, 37).Value = 29  # AK4
$ws.Cells.Item(4, 38).Value = 41  # AL4
$ws.Cells.Item(4, 39).Value = 451  # AM4

# --- Row 5 ---
$ws.Cells.Item(5, 33).Value = 8  # AG5
$ws.Cells.Item(5, 34).Value = 15  # AH5
$ws.Cells.Item(5, 35).Value = 12  # AI5
$ws.Cells.Item(5, 36).Value = 34  # AJ5
$ws.Cells.Item(5, 37).Value = 29  # AK5
$ws.Cells.Item(5, 38).Value = 41  # AL5
$ws.Cells.Item(5, 39).Value = 401  # AM5

# --- Row 6 ---
$ws.Cells.Item(6, 7).Value = 2.75  # G6
$ws.Cells.Item(6, 9).Value = 2.7  # I6
$ws.Cells.Item(6, 10).Value = 3.4  # J6
$ws.Cells.Item(6, 11).Value = 2  # K6
$ws.Cells.Item(6, 13).Value = 1.08  # M6
$ws.Cells.Item(6, 14).Value = 7.5  # N6
$ws.Cells.Item(6, 15).Value = 1.4  # O6
$ws.Cells.Item(6, 16).Value = 2.75  # P6
$ws.Cells.Item(6, 17).Value = 2.25  # Q6
$ws.Cells.Item(6, 18).Value = 1.62  # R6
$ws.Cells.Item(6, 19).Value = 1.5  # S6
$ws.Cells.Item(6, 20).Value = 2.5  # T6
$ws.Cells.Item(6, 21).Value = 1.91  # U6
$ws.Cells.Item(6, 22).Value = 1.8  # V6
$ws.Cells.Item(6, 23).Value = 7.5  # W6
$ws.Cells.Item(6, 27).Value = 23  # AA6
$ws.Cells.Item(6, 28).Value = 34  # AB6
$ws.Cells.Item(6, 29).Value = 7.5  # AC6
$ws.Cells.Item(6, 31).Value = 15  # AE6
$ws.Cells.Item(6, 32).Value = 51  # AF6
$ws.Cells.Item(6, 33).Value = 7.5  # AG6
$ws.Cells.Item(6, 34).Value = 13  # AH6
$ws.Cells.Item(6, 35).Value = 11  # AI6
$ws.Cells.Item(6, 36).Value = 29  # AJ6
$ws.Cells.Item(6, 37).Value = 23  # AK6
$ws.Cells.Item(6, 38).Value = 34  # AL6
$ws.Cells.Item(6, 39).Value = 351  # AM6
$ws.Cells.Item(6, 41).Value = 15  # AO6
$ws.Cells.Item(6, 45).Value = 201  # AS6
$ws.Cells.Item(6, 46).Value = 2.5  # AT6
$ws.Cells.Item(6, 50).Value = 15  # AX6
$ws.Cells.Item(6, 54).Value = 201  # BB6

# --- Row 7 ---
$ws.Cells.Item(7, 33).Value = 19  # AG7
$ws.Cells.Item(7, 34).Value = 51  # AH7
$ws.Cells.Item(7, 35).Value = 34  # AI7
$ws.Cells.Item(7, 36).Value = 151  # AJ7
$ws.Cells.Item(7, 37).Value = 81  # AK7
$ws.Cells.Item(7, 39).Value = 120  # AM7
